# Applies the "graphics design 2 added" edit to the document.
#
# Strategy: locate each target paragraph with Find (on stable anchor
# text), then overwrite that paragraph's Range with a minimal, exact
# OOXML fragment via Range.InsertXML. This lets us drop stray
# <w:lang>/<w:proofErr> markup and merge split runs precisely, instead
# of relying on Find&Replace (which only patches text and preserves
# whatever rPr happens to sit on the first matched run).

$d = $word.ActiveDocument

function New-WordOpenXml([string]$bodyXml) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Set-ParagraphXml([string]$findText, [string]$innerParagraphXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Anchor text not found: " + $findText)
    }
    $para = $rng.Paragraphs(1)
    $xml = New-WordOpenXml("<w:body>" + $innerParagraphXml + "</w:body>")
    $para.Range.InsertXML($xml) | Out-Null
}

# 1. "Inhoud website" (Kop1) - drop lang rPr + proofErr wrapping, merge runs
Set-ParagraphXml "Inhoud website" '<w:p><w:pPr><w:pStyle w:val="Kop1"/></w:pPr><w:r><w:t>Inhoud website</w:t></w:r></w:p>'

# 2. "Homepage " (Kop2) - drop lang rPr
Set-ParagraphXml "Homepage " '<w:p><w:pPr><w:pStyle w:val="Kop2"/></w:pPr><w:r><w:t xml:space="preserve">Homepage </w:t></w:r></w:p>'

# 3. "HEADER:" (bold) - drop lang rPr, keep bold
Set-ParagraphXml "HEADER:" '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>HEADER:</w:t></w:r></w:p>'

# 4. "Alle leerlingen..." paragraph -> "Wedstrijd 'Vergeten Geschiedenis'"
Set-ParagraphXml "Alle leerlingen uit de derde" ('<w:p><w:r><w:t>Wedstrijd ' + [char]0x2018 + 'Vergeten Geschiedenis' + [char]0x2019 + '</w:t></w:r></w:p>')

# 5. "International Convention Center" - drop fr-FR lang rPr (both on pPr and run)
Set-ParagraphXml "International Convention Center" '<w:p><w:r><w:t>International Convention Center (ICC) te Gent</w:t></w:r></w:p>'

# 6. "Adres<nbsp>: " paragraph - merge "Familie van " + "Rysselbeghedreef"
#    (+proofErr) + " 2, 9000 G" into one run. The "Adres" / ":" are
#    separated by a non-breaking space (U+00A0) in the source - keep it
#    untouched, it was not part of the diff.
$nbsp = [char]0x00A0
Set-ParagraphXml "Adres" ('<w:p><w:r><w:t xml:space="preserve">Adres' + $nbsp + ': </w:t></w:r><w:r><w:t>Familie van Rysselbeghedreef 2, 9000 G</w:t></w:r><w:r><w:t xml:space="preserve">ent </w:t></w:r></w:p>')

# 7. Insert new paragraph right after the first blank paragraph that
#    follows the "Geschiedenis " (Kop2) heading near the end of the
#    document. "Geschiedenis " also occurs earlier (inside "Het
#    Geschiedenis Olympiade "), so anchor on the unique "Deelnemen "
#    paragraph just before it and walk forward with .Next() instead of
#    a plain text Find.
$rng = $d.Content
$found = $rng.Find.Execute("Deelnemen ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Anchor paragraph 'Deelnemen ' not found"
}
$deelnemenPara = $rng.Paragraphs(1)
$headingPara = $deelnemenPara.Next().Next().Next()
if ($headingPara.Range.Text -notmatch "Geschiedenis") {
    throw ("Expected 'Geschiedenis ' heading paragraph, got: " + $headingPara.Range.Text)
}
$blankPara = $headingPara.Next()
$blankPara.Range.InsertParagraphAfter() | Out-Null
$newPara = $blankPara.Next()
$newText = "De Geschiedenis Olympiade " + [char]0x201C + "Vergeten Geschiedenissen" + [char]0x201D + " is een project van de Vakgroep Geschiedenis van de Universiteit Gent."
$newXml = New-WordOpenXml("<w:body><w:p><w:r><w:t>" + $newText + "</w:t></w:r></w:p></w:body>")
$newPara.Range.InsertXML($newXml) | Out-Null
